$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45624
$ws.Range("B2").Value = 45625
$ws.Range("C2").Value = 45626
$ws.Range("D2").Value = 45627
$ws.Range("E2").Value = 45628
